$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: USN, Name, Age, Phone, Vaccine_Dose
$ws.Range("A1").Value = "USN"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Age"
$ws.Range("D1").Value = "Phone"
$ws.Range("E1").Value = "Vaccine_Dose"

# Widen column E (Vaccine_Dose) to fit its header text
$ws.Columns.Item(5).ColumnWidth = 13

# Leave the cursor parked on F3, matching the author's saved selection
$ws.Range("F3").Select() | Out-Null
